# Edit script: merges the "Acompanha e altera..." runs, drops the stray
# _GoBack bookmark from the title, and appends the new "User Stories" section.

$d = $word.ActiveDocument

# 1) Merge the three runs that make up the "Acompanha e altera..." bullet
#    into a single run (Find/Replace naturally collapses the match into one
#    run since there is no formatting difference between the pieces).
$find = $d.Content.Find
$find.ClearFormatting()
$oldText = "Acompanha e altera andamento dos pedidos pelo painel (Menos os “Aguardando Aprovação”)"
$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $oldText, 2) | Out-Null

# 2) Drop the "_GoBack" bookmark that currently sits in the title paragraph.
#    (It gets re-added later, around the new "Encarregado" user story.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3) Append the new "Users Stories" section (4 blank paragraphs, the
#    "Users Stories" heading, and the four role paragraphs) right after the
#    last paragraph in the body, preserving the existing sectPr.
$fragment = '<w:p/><w:p/><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Users</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>tories</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Gerente – </w:t></w:r><w:r><w:t>Como gerente ele irá</w:t></w:r><w:r><w:t xml:space="preserve"> aprovar os </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>pedidos</w:t></w:r><w:r><w:t xml:space="preserve"> abertos </w:t></w:r><w:r><w:t>pelos vendedores</w:t></w:r><w:r><w:t xml:space="preserve"> e usuários alimentando seus prazos de entrega, </w:t></w:r><w:r><w:t>poderá também acompanha-los e realizar a</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>venda</w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>faturar) cada pedido. Na parte financeir</w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> poderá acompanhar seus saldos bancários e visões financeiras de contas a pagar e a receber. </w:t></w:r><w:r><w:t xml:space="preserve"> Extração de relatórios de centro de custo, vendas por gerente e venda por clientes.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Usuário Financeiro –</w:t></w:r><w:r><w:t xml:space="preserve"> Como usuário financeiro poderá cadastrar novos </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Clientes/Fornecedores</w:t></w:r><w:r><w:t xml:space="preserve"> e </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Produtos</w:t></w:r><w:r><w:t xml:space="preserve">, irá receber os pedidos feitos pelos clientes e lança-los no sistema, ficando com o estado de “Aguardando Aprovação”. Irá realizar a manutenção de </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>fluxo</w:t></w:r><w:r><w:t xml:space="preserve"> e conciliação bancaria.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Encarregado</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:t xml:space="preserve">Como encarregado poderá acompanhar e alterar os </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>estados dos pedidos</w:t></w:r><w:r><w:t xml:space="preserve"> na esteira, até estarem totalmente prontos para entrega.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Vendedor </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">– </w:t></w:r><w:r><w:t>Como Vendedor poderá apenas abrir pedidos, acompanhar a esteira de produção dos seus pedidos e acompanhar diariamente suas comissões.</w:t></w:r></w:p>'
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$end = $d.Content
$end.Collapse(0)
$end.InsertXML($xml)

Write-Output "done"
